$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 824.125
$ws.Range("I28").Value = 783.46155
$ws.Range("K28").Value = 783.46155
$ws.Range("M28").Value = -298.46155

$ws.Range("H41").Value = 384.26666
$ws.Range("I41").Value = 600.75
$ws.Range("J41").Value = 136.85715
$ws.Range("K41").Value = 600.75
$ws.Range("L41").Value = 136.85715
$ws.Range("M41").Value = -160.75
$ws.Range("N41").Value = -1016.85715

$ws.Range("H43").Value = 3042.7144
$ws.Range("J43").Value = 3859.8
$ws.Range("L43").Value = 3859.8
$ws.Range("N43").Value = -3997.8

$ws.Range("H80").Value = 742.25
$ws.Range("I80").Value = 606.3333
$ws.Range("J80").Value = 1150
$ws.Range("K80").Value = 1818.9999
$ws.Range("L80").Value = 3450
$ws.Range("M80").Value = -820.9999
$ws.Range("N80").Value = -5446

$ws.Range("H83").Value = 742.25
$ws.Range("I83").Value = 606.3333
$ws.Range("J83").Value = 1150
$ws.Range("K83").Value = 5456.9997
$ws.Range("L83").Value = 10350
$ws.Range("M83").Value = -464.9997000000003
$ws.Range("N83").Value = -20334

$ws.Range("H113").Value = 4882.9165
$ws.Range("I113").Value = 4698.125
$ws.Range("K113").Value = 4698.125
$ws.Range("M113").Value = -1444.125

$ws.Range("H137").Value = 44768.543
$ws.Range("I137").Value = 144638.28
$ws.Range("J137").Value = 3645.7058
$ws.Range("K137").Value = 433914.84
$ws.Range("L137").Value = 10937.1174
$ws.Range("M137").Value = -431364.84
$ws.Range("N137").Value = -16037.1174

$ws.Range("H138").Value = 2860.439
$ws.Range("I138").Value = 2674.3333
$ws.Range("J138").Value = 3055.85
$ws.Range("K138").Value = 8022.999899999999
$ws.Range("L138").Value = 9167.549999999999
$ws.Range("M138").Value = -2882.999899999999
$ws.Range("N138").Value = -19447.55

$ws.Range("H141").Value = 1954.5588
$ws.Range("I141").Value = 1954.5588
$ws.Range("K141").Value = 5863.6764
$ws.Range("M141").Value = -683.6764000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19376.213
$ws.Range("I32").Value = 19068.086
$ws.Range("K32").Value = 19068.086
$ws.Range("M32").Value = -18781.086

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H74").Value = 27358.275
$ws.Range("I74").Value = 27931.59
$ws.Range("K74").Value = 27931.59
$ws.Range("M74").Value = -27057.59

$ws.Range("H77").Value = 27358.275
$ws.Range("I77").Value = 27931.59
$ws.Range("K77").Value = 139657.95
$ws.Range("M77").Value = -135289.95

$ws.Range("H122").Value = 3392.6956
$ws.Range("I122").Value = 3335.8096
$ws.Range("K122").Value = 10007.4288
$ws.Range("M122").Value = -7557.4288

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2542.4092
$ws.Range("I20").Value = 1875.5
$ws.Range("J20").Value = 3342.7
$ws.Range("K20").Value = 1875.5
$ws.Range("L20").Value = 3342.7
$ws.Range("M20").Value = -1628.5
$ws.Range("N20").Value = -3836.7

$ws.Range("H86").Value = 3034.889
$ws.Range("I86").Value = 2789.25
$ws.Range("K86").Value = 2789.25
$ws.Range("M86").Value = -1666.25

$ws.Range("H89").Value = 3034.889
$ws.Range("I89").Value = 2789.25
$ws.Range("K89").Value = 13946.25
$ws.Range("M89").Value = -8330.25

$ws.Range("H101").Value = 69998
$ws.Range("J101").Value = 69998
$ws.Range("L101").Value = 69998
$ws.Range("N101").Value = -76488

$ws.Range("H134").Value = 1824.0851
$ws.Range("I134").Value = 1647.2727
$ws.Range("K134").Value = 4941.8181
$ws.Range("M134").Value = -2406.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2268.4055
$ws.Range("I31").Value = 2092.1614
$ws.Range("J31").Value = 3179
$ws.Range("K31").Value = 2092.1614
$ws.Range("L31").Value = 3179
$ws.Range("M31").Value = -1797.1614
$ws.Range("N31").Value = -3769

$ws.Range("H34").Value = 2268.4055
$ws.Range("I34").Value = 2092.1614
$ws.Range("J34").Value = 3179
$ws.Range("K34").Value = 2092.1614
$ws.Range("L34").Value = 3179
$ws.Range("M34").Value = -1890.1614
$ws.Range("N34").Value = -3583

$ws.Range("H86").Value = 2785.6072
$ws.Range("I86").Value = 2723.4092
$ws.Range("K86").Value = 2723.4092
$ws.Range("M86").Value = -1600.4092

$ws.Range("H89").Value = 2785.6072
$ws.Range("I89").Value = 2723.4092
$ws.Range("K89").Value = 13617.046
$ws.Range("M89").Value = -8001.046

$ws.Range("H122").Value = 2848.9092
$ws.Range("I122").Value = 2141.5
$ws.Range("J122").Value = 3697.8
$ws.Range("K122").Value = 6424.5
$ws.Range("L122").Value = 11093.4
$ws.Range("M122").Value = -3974.5
$ws.Range("N122").Value = -15993.4

$ws.Range("H132").Value = 2843.1482
$ws.Range("I132").Value = 2448.4092
$ws.Range("K132").Value = 7345.2276
$ws.Range("M132").Value = -4815.2276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 768.86365
$ws.Range("J113").Value = 742.5625
$ws.Range("L113").Value = 2227.6875
$ws.Range("N113").Value = -6567.6875

$ws.Range("H132").Value = 1211.75
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 1282.3334
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 11541.0006
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -16601.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2729.9565
$ws.Range("I122").Value = 2339.45
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 7018.349999999999
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -4568.349999999999
$ws.Range("N122").Value = -20900.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 46631.89
$ws.Range("I40").Value = 50633.184
$ws.Range("K40").Value = 50633.184
$ws.Range("M40").Value = -50497.184

$ws.Range("H46").Value = 13647.053
$ws.Range("I46").Value = 15819.6
$ws.Range("J46").Value = 5500
$ws.Range("K46").Value = 15819.6
$ws.Range("L46").Value = 5500
$ws.Range("M46").Value = -15631.6
$ws.Range("N46").Value = -5876

$ws.Range("H55").Value = 1562.5
$ws.Range("J55").Value = 1776.4
$ws.Range("L55").Value = 1776.4
$ws.Range("N55").Value = -2122.4

$ws.Range("H104").Value = 15550.286
$ws.Range("J104").Value = 15550.286
$ws.Range("L104").Value = 15550.286
$ws.Range("N104").Value = -22538.286

$ws.Range("H122").Value = 391814.2
$ws.Range("I122").Value = 10306.182
$ws.Range("K122").Value = 30918.546
$ws.Range("M122").Value = -28468.546

$ws.Range("H136").Value = 2900
$ws.Range("I136").Value = 2756.4102
$ws.Range("K136").Value = 8269.230599999999
$ws.Range("M136").Value = -5719.230599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 99648.086
$ws.Range("J62").Value = 146747.38
$ws.Range("L62").Value = 146747.38
$ws.Range("N62").Value = -147995.38

$ws.Range("H65").Value = 99648.086
$ws.Range("J65").Value = 146747.38
$ws.Range("L65").Value = 733736.9
$ws.Range("N65").Value = -739976.9

$ws.Range("H126").Value = 86672.21000000001
$ws.Range("I126").Value = 97867.95
$ws.Range("K126").Value = 293603.85
$ws.Range("M126").Value = -291133.85

$ws.Range("H132").Value = 60634.168
$ws.Range("I132").Value = 64024.53
$ws.Range("K132").Value = 192073.59
$ws.Range("M132").Value = -189543.59

$ws.Range("H136").Value = 4099.9033
$ws.Range("I136").Value = 3485.6072
$ws.Range("K136").Value = 10456.8216
$ws.Range("M136").Value = -7906.821599999999
